$wb = $excel.ActiveWorkbook

# The "SignUpNormalUser" sheet is duplicated (exactly like using Excel's
# "Move or Copy... > Create a copy" on its tab) to create the new
# "SignUpReseller" sheet right after it. Using .Copy() (rather than
# copy/paste of a range) preserves formulas, styles and number formats
# verbatim, and naturally makes the freshly created sheet the active tab -
# matching the target workbook exactly.
$src = $wb.Worksheets.Item("SignUpNormalUser")
$src.Copy([System.Reflection.Missing]::Value, $src)

$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "SignUpReseller"

$newSheet.Activate()
